$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 309.5
$ws.Range("I9").Value = 118.333336
$ws.Range("K9").Value = 118.333336
$ws.Range("M9").Value = 50.666664
$ws.Range("H17").Value = 5367.567
$ws.Range("J17").Value = 5487.1724
$ws.Range("L17").Value = 16461.5172
$ws.Range("N17").Value = -16797.5172
$ws.Range("H40").Value = 1098.3214
$ws.Range("I40").Value = 786.7059
$ws.Range("J40").Value = 1579.909
$ws.Range("K40").Value = 786.7059
$ws.Range("L40").Value = 1579.909
$ws.Range("M40").Value = -611.7059
$ws.Range("N40").Value = -1929.909
$ws.Range("H62").Value = 3675.8
$ws.Range("I62").Value = 3025.3333
$ws.Range("K62").Value = 3025.3333
$ws.Range("M62").Value = -2401.3333
$ws.Range("H65").Value = 3675.8
$ws.Range("I65").Value = 3025.3333
$ws.Range("K65").Value = 15126.6665
$ws.Range("M65").Value = -12006.6665
$ws.Range("H74").Value = 8932770
$ws.Range("J74").Value = 20837334
$ws.Range("L74").Value = 20837334
$ws.Range("N74").Value = -20839206
$ws.Range("H77").Value = 8932770
$ws.Range("J77").Value = 20837334
$ws.Range("L77").Value = 104186670
$ws.Range("N77").Value = -104196030
$ws.Range("H129").Value = 1400.7587
$ws.Range("J129").Value = 1539.4231
$ws.Range("L129").Value = 4618.2693
$ws.Range("N129").Value = -14618.2693
$ws.Range("H132").Value = 2250.8572
$ws.Range("I132").Value = 2250.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6752.571599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -4222.571599999999
$ws.Range("H138").Value = 129363.71
$ws.Range("I138").Value = 3414.1428
$ws.Range("J138").Value = 141441.06
$ws.Range("K138").Value = 10242.4284
$ws.Range("L138").Value = 424323.18
$ws.Range("M138").Value = -5102.428400000001
$ws.Range("N138").Value = -434603.18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2185.4827
$ws.Range("I45").Value = 2032.0555
$ws.Range("J45").Value = 2436.5454
$ws.Range("K45").Value = 2032.0555
$ws.Range("L45").Value = 2436.5454
$ws.Range("M45").Value = -1655.0555
$ws.Range("N45").Value = -3190.5454
$ws.Range("H61").Value = 1933.919
$ws.Range("I61").Value = 1636.9412
$ws.Range("J61").Value = 5299.6665
$ws.Range("K61").Value = 1636.9412
$ws.Range("L61").Value = 5299.6665
$ws.Range("M61").Value = -1424.9412
$ws.Range("N61").Value = -5723.6665
$ws.Range("H74").Value = 30304710
$ws.Range("I74").Value = 41667340
$ws.Range("J74").Value = 4357.1113
$ws.Range("K74").Value = 41667340
$ws.Range("L74").Value = 4357.1113
$ws.Range("M74").Value = -41666466
$ws.Range("N74").Value = -6105.1113
$ws.Range("H77").Value = 30304710
$ws.Range("I77").Value = 41667340
$ws.Range("J77").Value = 4357.1113
$ws.Range("K77").Value = 208336700
$ws.Range("L77").Value = 21785.5565
$ws.Range("M77").Value = -208332332
$ws.Range("N77").Value = -30521.5565
$ws.Range("H132").Value = 9900.661
$ws.Range("I132").Value = 1657.8864
$ws.Range("J132").Value = 30049.666
$ws.Range("K132").Value = 4973.6592
$ws.Range("L132").Value = 90148.99800000001
$ws.Range("M132").Value = -2443.6592
$ws.Range("N132").Value = -95208.99800000001
$ws.Range("H136").Value = 1933.919
$ws.Range("I136").Value = 1636.9412
$ws.Range("J136").Value = 5299.6665
$ws.Range("K136").Value = 4910.8236
$ws.Range("L136").Value = 15898.9995
$ws.Range("M136").Value = -2360.8236
$ws.Range("N136").Value = -20998.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2507.4075
$ws.Range("I20").Value = 2106.25
$ws.Range("J20").Value = 3090.9092
$ws.Range("K20").Value = 2106.25
$ws.Range("L20").Value = 3090.9092
$ws.Range("M20").Value = -1859.25
$ws.Range("N20").Value = -3584.9092
$ws.Range("H94").Value = 889.1111
$ws.Range("I94").Value = 647.7059
$ws.Range("K94").Value = 647.7059
$ws.Range("M94").Value = -196.7059
$ws.Range("H134").Value = 3203.9575
$ws.Range("I134").Value = 3165.348
$ws.Range("J134").Value = 4980
$ws.Range("K134").Value = 9496.044
$ws.Range("L134").Value = 14940
$ws.Range("M134").Value = -6961.044
$ws.Range("N134").Value = -20010
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 24415.908
$ws.Range("I58").Value = 1743.5834
$ws.Range("J58").Value = 51622.7
$ws.Range("K58").Value = 1743.5834
$ws.Range("L58").Value = 51622.7
$ws.Range("M58").Value = -1540.5834
$ws.Range("N58").Value = -52028.7
$ws.Range("H99").Value = 20003700
$ws.Range("I99").Value = 3182.9092
$ws.Range("J99").Value = 35718390
$ws.Range("K99").Value = 3182.9092
$ws.Range("L99").Value = 35718390
$ws.Range("M99").Value = -1684.9092
$ws.Range("N99").Value = -35721386
$ws.Range("H122").Value = 859.26666
$ws.Range("I122").Value = 859.26666
$ws.Range("K122").Value = 2577.79998
$ws.Range("M122").Value = -127.7999799999998
$ws.Range("H126").Value = 20003700
$ws.Range("I126").Value = 3182.9092
$ws.Range("J126").Value = 35718390
$ws.Range("K126").Value = 9548.7276
$ws.Range("L126").Value = 107155170
$ws.Range("M126").Value = -7078.7276
$ws.Range("N126").Value = -107160110
$ws.Range("H136").Value = 24415.908
$ws.Range("I136").Value = 1743.5834
$ws.Range("J136").Value = 51622.7
$ws.Range("K136").Value = 5230.7502
$ws.Range("L136").Value = 154868.1
$ws.Range("M136").Value = -2680.7502
$ws.Range("N136").Value = -159968.1
$ws.Range("H141").Value = 31332.104
$ws.Range("J141").Value = 31332.104
$ws.Range("L141").Value = 31332.104
$ws.Range("N141").Value = -41692.104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.111115
$ws.Range("I2").Value = 35
$ws.Range("K2").Value = 210
$ws.Range("M2").Value = -97
$ws.Range("H4").Value = 128
$ws.Range("I4").Value = 128
$ws.Range("K4").Value = 384
$ws.Range("M4").Value = -272
$ws.Range("H58").Value = 2400
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 3800
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 11400
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -11656
$ws.Range("H113").Value = 914.2857
$ws.Range("J113").Value = 930.7692
$ws.Range("L113").Value = 2792.3076
$ws.Range("N113").Value = -7132.3076
$ws.Range("H131").Value = 643.04
$ws.Range("J131").Value = 766.2465999999999
$ws.Range("L131").Value = 2298.7398
$ws.Range("N131").Value = -12378.7398
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 33468.668
$ws.Range("J3").Value = 100000
$ws.Range("L3").Value = 100000
$ws.Range("N3").Value = -100232
$ws.Range("H5").Value = 3666.6667
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224
$ws.Range("H70").Value = 2985007
$ws.Range("J70").Value = 5220429
$ws.Range("L70").Value = 5220429
$ws.Range("N70").Value = -5220969
$ws.Range("H73").Value = 2985007
$ws.Range("J73").Value = 5220429
$ws.Range("L73").Value = 5220429
$ws.Range("N73").Value = -5222301
$ws.Range("H102").Value = 2662
$ws.Range("I102").Value = 2382.6667
$ws.Range("K102").Value = 2382.6667
$ws.Range("M102").Value = -760.6667000000002
$ws.Range("H122").Value = 6937.615
$ws.Range("I122").Value = 7812.857
$ws.Range("K122").Value = 23438.571
$ws.Range("M122").Value = -20988.571
$ws.Range("H132").Value = 13237.66
$ws.Range("I132").Value = 3637.8147
$ws.Range("J132").Value = 24507.043
$ws.Range("K132").Value = 10913.4441
$ws.Range("L132").Value = 73521.129
$ws.Range("M132").Value = -8383.444100000001
$ws.Range("N132").Value = -78581.129
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4500.3335
$ws.Range("J22").Value = 1750
$ws.Range("L22").Value = 1750
$ws.Range("N22").Value = -2340
$ws.Range("H27").Value = 4500.3335
$ws.Range("J27").Value = 1750
$ws.Range("L27").Value = 1750
$ws.Range("N27").Value = -1964
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3312
$ws.Range("H55").Value = 141.75
$ws.Range("I55").Value = 103.8
$ws.Range("K55").Value = 103.8
$ws.Range("M55").Value = 69.2
$ws.Range("H119").Value = 29000
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676
$ws.Range("H122").Value = 1311541.9
$ws.Range("I122").Value = 1786775.2
$ws.Range("K122").Value = 5360325.6
$ws.Range("M122").Value = -5357875.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H136").Value = 32261304
$ws.Range("I136").Value = 43012330
$ws.Range("J136").Value = 8237.5
$ws.Range("K136").Value = 129036990
$ws.Range("L136").Value = 24712.5
$ws.Range("M136").Value = -129034440
$ws.Range("N136").Value = -29812.5
